# The document had accumulated SharePoint "document library" metadata:
# the Content-Type schema, the document-library form template, and the
# SharePoint "properties" custom XML parts that SharePoint stamps into
# a .docx while it lives in a document library (customXml/item1.xml,
# item2.xml, item3.xml and their itemProps*.xml siblings). Bringing the
# file up to date means stripping that stale, SharePoint-only custom
# XML back out of the package so the document is a plain, portable
# .docx again - the theme and every other part are left untouched.
#
# Word exposes those package-level parts through
# Document.CustomXMLParts. Delete every part that isn't one of Word's
# own built-in parts (core/extended properties, cover-page props, ...),
# which aren't real package parts and can't be removed. Walk back to
# front since indices shift after each delete.

$d = $word.ActiveDocument

function Remove-CustomXmlPart($part) {
    $isBuiltIn = $false
    try {
        $isBuiltIn = [bool]$part.BuiltIn
    } catch {
        $isBuiltIn = $false
    }

    if ($isBuiltIn) {
        return
    }

    try {
        $part.Delete()
    } catch {
        # Could not drop the part through this call shape - caller
        # tries the next strategy.
    }
}

# Primary approach: walk the live collection back-to-front.
$parts = $d.CustomXMLParts
$count = 0
try { $count = [int]$parts.Count } catch { $count = 0 }

for ($i = $count; $i -ge 1; $i--) {
    try {
        $part = $parts.Item($i)
        Remove-CustomXmlPart $part
    } catch {
        # keep going - don't let one bad index stop the cleanup
    }
}

# Belt-and-braces: the SharePoint parts are also reachable by the
# namespaces they declare, so sweep those explicitly too in case the
# host's collection indices don't line up with a simple Count/Item
# walk.
$sharePointNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

foreach ($ns in $sharePointNamespaces) {
    try {
        $scoped = $d.CustomXMLParts.SelectByNamespace($ns)
        $scopedCount = 0
        try { $scopedCount = [int]$scoped.Count } catch { $scopedCount = 0 }
        for ($j = $scopedCount; $j -ge 1; $j--) {
            try {
                $scoped.Item($j).Delete()
            } catch {
            }
        }
    } catch {
        # Namespace not present (already removed) - nothing to do.
    }
}

try {
    Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
} catch {
    Write-Output "CustomXMLParts cleanup attempted"
}
